# Remove the existing slide comment (feedback has been addressed / tweaked away)
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = $s.Comments.Count; $i -ge 1; $i--) {
    $s.Comments.Item($i).Delete()
}
